# "making patterns for ulna"
# Insert a new "axis" worksheet between "trait" and "structures", populate it
# with Trait/Axis/Structure/Pattern-name data, add a new "in oba" cell on the
# trait sheet, and update the active selections accordingly.

$wb = $excel.ActiveWorkbook
$trait = $wb.Worksheets.Item("trait")
$structures = $wb.Worksheets.Item("structures")

# --- new "axis" sheet, placed right after "trait" ---------------------------
$axis = $wb.Worksheets.Add($null, $trait)
$axis.Name = "axis"

# header row
$axis.Range("A1").Value = "Trait"
$axis.Range("B1").Value = "Axis"
$axis.Range("C1").Value = "Structure"
$axis.Range("D1").Value = "Pattern name"

# row 2 is blank but carries the same "bold/explicit-font" style as column C
# of the trait sheet (and column A below) -- copy formats only.
$trait.Range("C2").Copy() | Out-Null
$axis.Range("A2").PasteSpecial(-4122) | Out-Null

# row 3 - ulna length of the olecranon
$trait.Range("C3").Copy() | Out-Null
$axis.Range("A3").PasteSpecial(-4122) | Out-Null
$axis.Range("A3").Value = "ulna length of the olecranon"
$axis.Range("B3").Value = "proximal-distal"
$axis.Range("C3").Value = "olecranon"

# row 4 - ulna proximal articular breadth
$trait.Range("C4").Copy() | Out-Null
$axis.Range("A4").PasteSpecial(-4122) | Out-Null
$axis.Range("A4").Value = "ulna proximal articular breadth"
$axis.Range("B4").Value = "medial-lateral"
$axis.Range("C4").Value = "humeral facet on radius"

# row 5 - ulna smallest depth of the olecranon
$trait.Range("C5").Copy() | Out-Null
$axis.Range("A5").PasteSpecial(-4122) | Out-Null
$axis.Range("A5").Value = "ulna smallest depth of the olecranon"
$axis.Range("B5").Value = "anterior-posterior"
$axis.Range("C5").Value = "'distalmost part of' some olecranon"

# row 6 - ulna depth across the process anaconaeus
$trait.Range("C6").Copy() | Out-Null
$axis.Range("A6").PasteSpecial(-4122) | Out-Null
$axis.Range("A6").Value = "ulna depth across the process anaconaeus"
$axis.Range("B6").Value = "anterior-posterior"
$axis.Range("C6").Value = "'proximalmost part of' some olecranon"

# row 7 - ulna lateral length
$trait.Range("C7").Copy() | Out-Null
$axis.Range("A7").PasteSpecial(-4122) | Out-Null
$axis.Range("A7").Value = "ulna lateral length"
$axis.Range("B7").Value = "proximal-distal"
$axis.Range("C7").Value = "'lateral side of' some ulna"

# size column A to its (longest) content, like a manual double-click autofit
$axis.Columns("A").AutoFit() | Out-Null

# --- trait sheet: new synonym/axis-header cell -------------------------------
$trait.Range("A2").Value = "in oba"

# --- selections / active cells ----------------------------------------------
$axis.Range("B8").Select() | Out-Null
$trait.Range("C3").Select() | Out-Null
$trait.Activate() | Out-Null
